$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1145
$ws1.Range("F3").Value = 0
$ws1.Range("F5").Value = 4942
$ws1.Range("F7").Value = 9103
$ws1.Range("F8").Value = 0
$ws1.Range("F10").Value = 0
$ws1.Range("F11").Value = 0
$ws1.Range("F12").Value = 0

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 0
$ws2.Range("F5").Value = 0
$ws2.Range("F6").Value = 0

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1145
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 0
$ws4.Range("F5").Value = 0
$ws4.Range("F6").Value = 0
$ws4.Range("F7").Value = 4942
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 9103
$ws4.Range("F11").Value = 0
$ws4.Range("F12").Value = 0
$ws4.Range("F14").Value = 0
$ws4.Range("F16").Value = 623
$ws4.Range("F17").Value = 0
